# Singapore Premier League workbook update (2024-06-14 20:31 refresh)
#
# The underlying data refresh did two things to "Singapore Premier League.xlsx":
#  1. Swapped the two team names "Albirex Niigata Singapore" and "Young Lions"
#     wherever they appear as HomeTeam (col E) / AwayTeam (col F) values.
#  2. Re-ordered several same-day fixture pairs, so that the full match
#     record (id number, teams, score, odds, ...) moved from one row to the
#     other while the running row index (col A) stayed put.
#
# Both effects are reproduced below by operating directly on the cell values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap whole match records between these row pairs (every column except
#    the running index "A", and the constant "Div"/"Date" columns C/D, which
#    are identical for both rows in a pair anyway).
# ---------------------------------------------------------------------------

$swapCols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

function Swap-MatchRows($rowA, $rowB) {
    foreach ($col in $swapCols) {
        $refA = $ws.Range("$col$rowA")
        $refB = $ws.Range("$col$rowB")
        $valA = $refA.Value2
        $valB = $refB.Value2
        $refA.Value = $valB
        $refB.Value = $valA
    }
}

Swap-MatchRows 4 5
Swap-MatchRows 18 19
Swap-MatchRows 26 27
Swap-MatchRows 36 37
Swap-MatchRows 38 39
Swap-MatchRows 56 57

# ---------------------------------------------------------------------------
# 2) For every other row, flip "Albirex Niigata Singapore" <-> "Young Lions"
#    wherever it shows up as HomeTeam (E) or AwayTeam (F). (Rows already
#    handled by the full swap above already have the correct team names.)
# ---------------------------------------------------------------------------

function Flip-TeamName($rowNum) {
    foreach ($col in @("E","F")) {
        $rng = $ws.Range("$col$rowNum")
        $v = $rng.Value2
        if ($v -eq "Albirex Niigata Singapore") {
            $rng.Value = "Young Lions"
        } elseif ($v -eq "Young Lions") {
            $rng.Value = "Albirex Niigata Singapore"
        }
    }
}

$singleRows = @(3,9,11,13,14,16,24,25,28,31,33,40,42,45,46,48,51,53,60,61,63,64,67,68)
foreach ($r in $singleRows) {
    Flip-TeamName $r
}
